$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "score" column (B) for rows 4-27 was stored as text, which breaks any
# SUM/aggregate formulas over the range. Fix it by converting those cells to
# real numeric values (the displayed numbers stay the same).
for ($r = 4; $r -le 27; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $cell.Value2 = [double]$cell.Value2
}

# Weekly update: append this week's two results.
$ws.Cells.Item(28, 1).Value2 = "Philadelphia Eagles"
$ws.Cells.Item(28, 2).NumberFormat = "@"
$ws.Cells.Item(28, 2).Value2 = "21"

$ws.Cells.Item(29, 1).Value2 = "Kansas City Chiefs"
$ws.Cells.Item(29, 2).NumberFormat = "@"
$ws.Cells.Item(29, 2).Value2 = "17"
